# The workbook lists two worker records in the "Estado de Cuenta" table
# (rows 16-17). This edit swaps the order in which the two workers are
# listed: the record that used to be on row 16 (BETILDA MUÑOZ BELTRAN)
# moves to row 17, and the record that used to be on row 17
# (MABEL POLO MENDOZA) moves to row 16. Cell formatting/styles stay with
# their row; only the data values move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for the two data rows (columns C-G).
$row16 = @{
    C = $ws.Range("C16").Value2
    D = $ws.Range("D16").Value2
    E = $ws.Range("E16").Value2
    F = $ws.Range("F16").Value2
    G = $ws.Range("G16").Value2
}

$row17 = @{
    C = $ws.Range("C17").Value2
    D = $ws.Range("D17").Value2
    E = $ws.Range("E17").Value2
    F = $ws.Range("F17").Value2
    G = $ws.Range("G17").Value2
}

# Write row 17's former data into row 16.
$ws.Range("C16").Value2 = $row17.C
$ws.Range("D16").Value2 = $row17.D
$ws.Range("E16").Value2 = $row17.E
$ws.Range("F16").Value2 = $row17.F
$ws.Range("G16").Value2 = $row17.G

# Write row 16's former data into row 17.
$ws.Range("C17").Value2 = $row16.C
$ws.Range("D17").Value2 = $row16.D
$ws.Range("E17").Value2 = $row16.E
$ws.Range("F17").Value2 = $row16.F
$ws.Range("G17").Value2 = $row16.G
